$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the description text for the "Correlation between poles" row (column F, row 3)
$ws.Range("F3").Value = "0: no correlation. 1: perfect correlation /!\ Feature not implement in the current version (0 is mandatory)"

# Make that cell's font red to flag the note
$ws.Range("F3").Font.Color = 255

# Update the selected cell to match the saved state
$ws.Range("F3").Select()
